# Insert a new weekly price record for Cebollín (Feria Lagunitas de Puerto Montt)
# as row 262, pushing the existing rows 262-397 down to 263-398.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(262).Insert()

$ws.Range("A262").Value = 4
$ws.Range("B262").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C262").Value = "Los Lagos"
$ws.Range("D262").Value = 44960
$ws.Range("E262").Value = 10
$ws.Range("F262").Value = 100112037
$ws.Range("G262").Value = "Cebollín"
$ws.Range("H262").Value = "Sin especificar"
$ws.Range("I262").Value = "Primera"
$ws.Range("J262").Value = 180
$ws.Range("K262").Value = 6000
$ws.Range("L262").Value = 6000
$ws.Range("M262").Value = 6000
$ws.Range("N262").Value = "$/paquete 36 unidades"
$ws.Range("O262").Value = "Región Metropolitana"
$ws.Range("P262").Value = 167
$ws.Range("Q262").Value = 36
$ws.Range("R262").Value = "Hortaliza"
